$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in missing grade values for row 8 (columns C-F) and row 15 (column F)
$ws.Range("C8").Value = 5
$ws.Range("D8").Value = 5
$ws.Range("E8").Value = 5
$ws.Range("F8").Value = 5

$ws.Range("F15").Value = 5

# Move the active selection to F15 (previously F31), matching the updated scroll position
$ws.Range("F15").Select()
